$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 18
$ws.Range("A19").HorizontalAlignment = -4108

$ws.Range("B19").Value = "Flatten Binary Tree to Linked List"
$ws.Range("B19").HorizontalAlignment = -4131

$ws.Range("D19").Value = "Tree"
$ws.Range("D19").HorizontalAlignment = -4108

$ws.Range("E19").Value = "medium"
$ws.Range("E19").HorizontalAlignment = -4108

$ws.Range("F19").Value = "leetcode 114"
$ws.Range("F19").HorizontalAlignment = -4108

$ws.Range("F23").Select()
